# Atualização agendada das bases de dados
# Updates a handful of numeric values on Sheet1 (column D) as scheduled
# in the upstream data refresh.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D3").Value  = 6.2
$ws.Range("D6").Value  = 12.1
$ws.Range("D26").Value = 10.5
$ws.Range("D27").Value = 8.699999999999999
$ws.Range("D30").Value = 9.1
$ws.Range("D33").Value = 13.7
$ws.Range("D35").Value = 15.1
$ws.Range("D38").Value = 12
$ws.Range("D39").Value = 11.3
